$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "PLO" (sheet1)
# ---------------------------------------------------------------
$plo = $wb.Worksheets.Item("PLO")

$plo.Range("B1").Value = "* ถ้า import ด้วย template นี้ จะต้องใส่ Sub-PLO ไปพร้อมกับ PLO"

$plo.Range("D2").Value = "Example:"
$plo.Range("E2").Value = "regular"
$plo.Range("F2").Value = "2564 (พ.ศ.)"

$plo.Range("B3").ClearContents()

$plo.Range("E4").Value = "Example:"
$plo.Range("F4").Value = 1
$plo.Range("G4").Value = "สามารถ…"
$plo.Range("H4").Value = "Able to… "

$plo.Range("A5:C5").ClearContents()
$plo.Range("H5").Value = "(Optional)"

$plo.Range("A6").ClearContents()

# ---------------------------------------------------------------
# Sheet "Sub-PLO" (sheet2)
# ---------------------------------------------------------------
$subplo = $wb.Worksheets.Item("Sub-PLO")

$subplo.Range("B1").Value = "* ถ้า import ด้วย template นี้ จะต้องใส่ Sub-PLO ไปพร้อมกับ PLO"

$subplo.Range("F2").Value = "Example:"
$subplo.Range("G2").Value = 1
$subplo.Range("H2").Value = 2
$subplo.Range("I2").Value = "พัฒนา..."
$subplo.Range("J2").Value = "Design…"

$subplo.Range("A3:D3").ClearContents()
$subplo.Range("G3").Value = "แปลว่า Sub PLO 1.2"
$subplo.Range("J3").Value = "(Optional)"

$subplo.Range("A4:B5").ClearContents()

# ---------------------------------------------------------------
# Sheet "PO" (sheet3)
# ---------------------------------------------------------------
$po = $wb.Worksheets.Item("PO")

$po.Range("B1").Value = "* ไม่จำเป็นต้อง import ไปพร้อมกับ PLO, Sub-PLO"

$po.Range("E2").Value = "Example:"
$po.Range("F2").Value = 1
$po.Range("G2").Value = "การ…"
$po.Range("H2").Value = "สามารถ…"

$po.Range("A3:C3").ClearContents()
$po.Range("A4").ClearContents()

# ---------------------------------------------------------------
# Selections (all sheets end up with B1 selected; PLO stays active tab)
# ---------------------------------------------------------------
$plo.Range("B1").Select() | Out-Null
$subplo.Range("B1").Select() | Out-Null
$po.Range("B1").Select() | Out-Null
$plo.Activate()
